$wb = $excel.ActiveWorkbook

# --- UDFData sheet: remove the now-unused "udf data N" placeholder entries ---
# (these were filler values in columns B, E, F, G, J, K of row 2; clearing them
#  drops the corresponding shared strings entirely on save)
$wsUDF = $wb.Worksheets.Item("UDFData")
$wsUDF.Range("B2").Clear()
$wsUDF.Range("E2").Clear()
$wsUDF.Range("F2").Clear()
$wsUDF.Range("G2").Clear()
$wsUDF.Range("J2").Clear()
$wsUDF.Range("K2").Clear()

# --- Update selections on the sheets that had their selection move ---
$wsName = $wb.Worksheets.Item("NameData")
$wsName.Range("E6").Select()

# --- Make UDFData the active/selected sheet with K2 selected ---
$wsUDF.Activate()
$wsUDF.Range("K2").Select()
